$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# AMD samples (rows 29-59) previously had no CPU time recorded (column K).
# Running them on CPU now populates column K with the measured times,
# stored as text to match the precision/format already used by column L (GPU time).
$kRange = $ws.Range("K29:K59")
$originalFormat = $kRange.NumberFormat

# Temporarily mark the range as Text so the values below are stored verbatim
# (preserving trailing zeros / exact formatting) instead of being parsed as numbers.
$kRange.NumberFormat = "@"

$ws.Cells.Item(29, 11).Value = "0.0528379000"
$ws.Cells.Item(30, 11).Value = "0.0140000000"
$ws.Cells.Item(31, 11).Value = "0.5789000000"
$ws.Cells.Item(32, 11).Value = "0.4280000000"
$ws.Cells.Item(33, 11).Value = "0.6647100000"
$ws.Cells.Item(34, 11).Value = "0.6452500000"
$ws.Cells.Item(35, 11).Value = "0.7257500000"
$ws.Cells.Item(36, 11).Value = "0.3880000000"
$ws.Cells.Item(37, 11).Value = "0.6890500000"
$ws.Cells.Item(38, 11).Value = "0.6660800000"
$ws.Cells.Item(39, 11).Value = "0.6112400000"
$ws.Cells.Item(40, 11).Value = "0.6583600000"
$ws.Cells.Item(41, 11).Value = "0.6252900000"
$ws.Cells.Item(42, 11).Value = "0.6474900000"
$ws.Cells.Item(43, 11).Value = "0.6949400000"
$ws.Cells.Item(44, 11).Value = "0.7830000000"
$ws.Cells.Item(45, 11).Value = "0.6440800000"
$ws.Cells.Item(46, 11).Value = "23.4670000000"
$ws.Cells.Item(47, 11).Value = "1.2809300000"
$ws.Cells.Item(48, 11).Value = "0.0140000000"
$ws.Cells.Item(49, 11).Value = "0.0050000000"
$ws.Cells.Item(50, 11).Value = "0.1110000000"
$ws.Cells.Item(51, 11).Value = "0.5719000000"
$ws.Cells.Item(52, 11).Value = "0.6080700000"
$ws.Cells.Item(53, 11).Value = "0.0130000000"
$ws.Cells.Item(54, 11).Value = "0.5607500000"
$ws.Cells.Item(55, 11).Value = "0.5937400000"
$ws.Cells.Item(56, 11).Value = "0.5751100000"
$ws.Cells.Item(57, 11).Value = "0.0380000000"
$ws.Cells.Item(58, 11).Value = "4.3045000000"
$ws.Cells.Item(59, 11).Value = "0.6114800000"

# Restore the original number format now that the cells hold text values;
# the engine reuses the existing style since the effective formatting matches.
$kRange.NumberFormat = $originalFormat

Write-Host "Populated CPU time (column K) for AMD sample rows 29-59"
